$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D/E hold plain text in the source data (prices use "." as a
# thousands separator, e.g. "69.040.31", and percentages keep their
# original padding, e.g. "  -0.21%  "). Cells whose new value parses as a
# clean decimal number (e.g. "629.37") would otherwise be auto-converted
# to a numeric type by Excel on assignment, so force those to Text first.
$preserveAsText = @(
    "D5", "D6", "D9", "D11", "D12", "D13", "D14", "D18",
    "D21", "D23", "D24", "D26", "D27", "D28", "D33", "D34",
    "D35", "D40", "D41", "D42", "D43", "D45", "D46", "D47",
    "D48", "D50", "D51"
)
foreach ($addr in $preserveAsText) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "69.040.31"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "3.771.91"
$ws.Range("E3").Value = "  -1.42%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "629.37"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("D6").Value = "165.22"
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("D7").Value = "3.771.13"
$ws.Range("E7").Value = "  -1.41%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "0.521"
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("E10").Value = "  -2.38%  "
$ws.Range("D11").Value = "0.454"
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D12").Value = "6.89"
$ws.Range("E12").Value = "  +4.59%  "
$ws.Range("D13").Value = "0.0000238"
$ws.Range("E13").Value = "  -5.12%  "
$ws.Range("D14").Value = "34.79"
$ws.Range("E14").Value = "  -3.71%  "
$ws.Range("D15").Value = "4.409.05"
$ws.Range("E15").Value = "  -1.38%  "
$ws.Range("D16").Value = "3.775.01"
$ws.Range("E16").Value = "  +0.42%  "
$ws.Range("D17").Value = "69.061.81"
$ws.Range("E17").Value = "  -0.26%  "
$ws.Range("D18").Value = "17.64"
$ws.Range("E18").Value = "  -2.75%  "
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("E20").Value = "  -2.15%  "
$ws.Range("D21").Value = "467.85"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("E22").Value = "  -1.99%  "
$ws.Range("D23").Value = "0.702"
$ws.Range("E23").Value = "  -1.04%  "
$ws.Range("D24").Value = "81.88"
$ws.Range("E24").Value = "  -2.48%  "
$ws.Range("E25").Value = "  -7.34%  "
$ws.Range("D26").Value = "12.08"
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("D27").Value = "2.11"
$ws.Range("E27").Value = "  -2.02%  "
$ws.Range("D28").Value = "10.13"
$ws.Range("E28").Value = "  +0.61%  "
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").Value = "3.922.55"
$ws.Range("E30").Value = "  -1.36%  "
$ws.Range("E31").Value = "  +2.23%  "
$ws.Range("E32").Value = "  -0.90%  "
$ws.Range("D33").Value = "7.10"
$ws.Range("E33").Value = "  -3.08%  "
$ws.Range("D34").Value = "0.177"
$ws.Range("E34").Value = "  +18.70%  "
$ws.Range("D35").Value = "28.40"
$ws.Range("E35").Value = "  -3.05%  "
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("D37").Value = "3.726.03"
$ws.Range("E37").Value = "  -1.26%  "
$ws.Range("E38").Value = "  -2.77%  "
$ws.Range("E39").Value = "  -0.87%  "
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").Value = "5.79"
$ws.Range("E40").Value = "  -2.39%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "3.24"
$ws.Range("E41").Value = "  -5.30%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "0.962"
$ws.Range("E43").Value = "  -2.28%  "
$ws.Range("D45").Value = "1.98"
$ws.Range("E45").Value = "  +4.55%  "
$ws.Range("D46").Value = "156.13"
$ws.Range("E46").Value = "  +0.81%  "
$ws.Range("D47").Value = "43.90"
$ws.Range("E47").Value = "  +3.20%  "
$ws.Range("D48").Value = "46.94"
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("E49").Value = "  -3.32%  "
$ws.Range("D50").Value = "0.293"
$ws.Range("E50").Value = "  -2.27%  "
$ws.Range("D51").Value = "8.34"
$ws.Range("E51").Value = "  -1.53%  "
